# edit.ps1
# Updates the DYNGROWTH holdings workbook:
#  - Refresh the "as of" date in the confidential disclosure note (A80)
#  - Refresh Weight (col D) and Percent Change (col E) values for rows 2-77

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; temporarily unprotect to make edits, then restore
$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

# --- Update the confidential disclosure text (cell A80) ---
$oldNote = $ws.Range("A80").Value
$newNote = $oldNote -replace "2021-04-22", "2021-04-23"
$ws.Range("A80").Value = $newNote

# --- Update Weight / Percent Change values for rows 2-77 ---
$data = @(
    @{ Row=2; D=0.06721300477059815; E=0.01803850234955284 },
    @{ Row=3; D=0.04013558134399115; E=0.009622126054686708 },
    @{ Row=4; D=0.03431156528678295; E=0.01547614418478038 },
    @{ Row=5; D=0.03036233974535597; E=0.01067803313266258 },
    @{ Row=6; D=0.02732097517375642; E=0.02104753786869806 },
    @{ Row=7; D=0.02383281605405788; E=0.01913550926240082 },
    @{ Row=8; D=0.1737044698626828; E=0.02596359743040688 },
    @{ Row=9; D=0.02470958617170105; E=0.002058360576341123 },
    @{ Row=10; D=0.02286111574014862; E=-0.005125157840006 },
    @{ Row=11; D=0.02270885546794332; E=0.01053268765133186 },
    @{ Row=12; D=0.02068928427378525; E=0.001422630772597966 },
    @{ Row=13; D=0.01892102743226344; E=0.02137643378519294 },
    @{ Row=14; D=0.01728190476743408; E=0.007961165048543606 },
    @{ Row=15; D=0.01761868864089934; E=0.0146220570012392 },
    @{ Row=16; D=0.01575808913827373; E=0.02298850574712641 },
    @{ Row=17; D=0.0146819170921875; E=0.001646738741319043 },
    @{ Row=18; D=0.01476411176382674; E=-0.001643047853768564 },
    @{ Row=19; D=0.01318721123360278; E=0.01554701200593556 },
    @{ Row=20; D=0.01251364348667226; E=0.005427899402930869 },
    @{ Row=21; D=0.01255213307221115; E=0.001275510204081565 },
    @{ Row=22; D=0.01290310795883202; E=0.01208541572012733 },
    @{ Row=23; D=0.01239946849763674; E=0.009651502504173681 },
    @{ Row=24; D=0.01282576491562206; E=0.009507237605285601 },
    @{ Row=25; D=0.01169703356155205; E=0.007559260872270324 },
    @{ Row=26; D=0.009589567032350341; E=0.01558257584700762 },
    @{ Row=27; D=0.009938399116414962; E=0.02107267224265308 },
    @{ Row=28; D=0.0106906441043735; E=0.01447685896029816 },
    @{ Row=29; D=0.01062999874900761; E=0.004202767360662163 },
    @{ Row=30; D=0.01020701761044896; E=0.009110354115503538 },
    @{ Row=31; D=0.009865907701634265; E=-0.05322039315966109 },
    @{ Row=32; D=0.01011726248450744; E=-0.0001838235294118196 },
    @{ Row=33; D=0.009494960277979549; E=0.0003491620111730764 },
    @{ Row=34; D=0.009244292809133857; E=0.0005510653930931309 },
    @{ Row=35; D=0.009136829239425494; E=0.01119518562768262 },
    @{ Row=36; D=0.008987641665225397; E=-0.0015384615384616 },
    @{ Row=37; D=0.008575253248737322; E=0.01110325318246108 },
    @{ Row=38; D=0.00872917116065596; E=0.01349192013227896 },
    @{ Row=39; D=0.008748456383662314; E=0.0005175983436851439 },
    @{ Row=40; D=0.007754721590636784; E=0.005943536404160454 },
    @{ Row=41; D=0.007360082048152461; E=0.02421392630353103 },
    @{ Row=42; D=0.007526412042802651; E=0.01850041362713406 },
    @{ Row=43; D=0.007870230777490356; E=0.001335648457325922 },
    @{ Row=44; D=0.007311889205755031; E=0.007818547762811434 },
    @{ Row=45; D=0.007830690005791794; E=0.008013052188100112 },
    @{ Row=46; D=0.007153240956117854; E=0.04015192620727093 },
    @{ Row=47; D=0.007564659046920074; E=-0.001026167265264277 },
    @{ Row=48; D=0.007126718720704504; E=0.01330897703549039 },
    @{ Row=49; D=0.007030009594014359; E=0.01414768806073163 },
    @{ Row=50; D=0.006775945985268176; E=-0.001706484641638251 },
    @{ Row=51; D=0.006477853848526355; E=0.005330071213246468 },
    @{ Row=52; D=0.006430995203946975; E=0.02244379620781567 },
    @{ Row=53; D=0.005415630234174277; E=0.0257558790593504 },
    @{ Row=54; D=0.006049980651301522; E=0.008928094092488514 },
    @{ Row=55; D=0.005756457131330598; E=0.03466076696165188 },
    @{ Row=56; D=0.005722107602051356; E=0.008334605144901586 },
    @{ Row=57; D=0.005741417083199856; E=0.008168553884288388 },
    @{ Row=58; D=0.005618185721096361; E=0.007253886010362809 },
    @{ Row=59; D=0.005080382709711619; E=0.0224896146683855 },
    @{ Row=60; D=0.005052809288138593; E=-0.008449622327486805 },
    @{ Row=61; D=0.004536272581368824; E=0.02362745098039221 },
    @{ Row=62; D=0.004478457342586672; E=0.02554843369143267 },
    @{ Row=63; D=0.004538860116531102; E=-0.05857621321171524 },
    @{ Row=64; D=0.004290133299057117; E=-0.004674306393244732 },
    @{ Row=65; D=0.004158330726728576; E=-0.004978026679111736 },
    @{ Row=66; D=0.003834727110496169; E=-0.004301619433198289 },
    @{ Row=67; D=0.003991192127340174; E=0.007354281893879433 },
    @{ Row=68; D=0.003354173314576831; E=0.004339336081579326 },
    @{ Row=69; D=0.003595541828933086; E=0.008950658930418687 },
    @{ Row=70; D=0.003040192094729145; E=0.0117027501462843 },
    @{ Row=71; D=0.00318109147036257; E=0.01433637091546869 },
    @{ Row=72; D=0.002401596502726288; E=0.02794565748051392 },
    @{ Row=73; D=0.002057009593537282; E=-0.006348519988993218 },
    @{ Row=74; D=0.002050985488237603; E=0.01685426481893981 },
    @{ Row=75; D=0.001523168745369778; E=0.01337792642140467 },
    @{ Row=76; D=0.001408751174912792; E=0.01084835265755957 },
    @{ Row=77; D=1; E=0.01248775196295471 }
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 4).Value = $row.D
    $ws.Cells.Item($row.Row, 5).Value = $row.E
}

# Restore protection state
if ($wasProtected) {
    $ws.Protect()
}
